# The two worksheets' contents were swapped: what used to live on Sheet1
# (Street/Town/State/LATITUDE/LONGITUDE, A1:E4) now lives on Sheet2, and
# what used to live on Sheet2 (Street/Town/State, A1:C4) now lives on
# Sheet1. Capture both ranges first (via Value2, which round-trips a real
# 2-D array through this COM shim, unlike Value), then clear and rewrite
# each sheet with the other's data so the swap is exact (including the
# used range / dimension).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Snapshot the full used ranges before mutating anything.
$sheet1Data = $ws1.Range("A1:E4").Value2
$sheet2Data = $ws2.Range("A1:C4").Value2

# Wipe both sheets clean so no stray columns/rows survive the swap.
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# Sheet1 now gets the old Sheet2 data (A1:C4).
$ws1.Range("A1:C4").Value2 = $sheet2Data

# Sheet2 now gets the old Sheet1 data (A1:E4).
$ws2.Range("A1:E4").Value2 = $sheet1Data
